# VNC Phidgets workbook edit
# - Insert "Name, IP" and "Serial #" columns after "Name"
# - Add "Notes" column at the end
# - Split the PhidgetSBC2 (1072) row into three rows, one per serial number,
#   and add the PSBC network name/IP for PhidgetSBC/PhidgetSBC2 rows
# - Re-wrap the "Description" column, keep "Printed Manual" centered+wrapped header
# - Resize/rebuild Table1 over the new A4:H31 range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Detach the existing table so we can freely restructure the grid
# ------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Unlist()

# ------------------------------------------------------------------
# 2. Make room: insert two blank columns (C,D) and two blank rows (14,15)
# ------------------------------------------------------------------
$ws.Range("C:D").Insert()
$ws.Range("14:15").Insert()

# ------------------------------------------------------------------
# 3. Write every cell of the final A4:H31 grid
# ------------------------------------------------------------------
$ws.Range("A4").Value = 'Model'
$ws.Range("B4").Value = 'Name'
$ws.Range("C4").Value = 'Name, IP'
$ws.Range("D4").Value = 'Serial #'
$ws.Range("E4").Value = 'Description'
$ws.Range("F4").Value = 'Count'
$ws.Range("G4").Value = 'Printed Manual'
$ws.Range("H4").Value = 'Notes'
$ws.Range("A5").Value = 1012
$ws.Range("B5").Value = 'PhidgetInterfaceKit 8/8/8'
$ws.Range("E5").Value = '8 Analog inputs, 8 Digital Inputs, 8 Digital Outputs'
$ws.Range("G5").Value = 'X'
$ws.Range("A6").Value = 1018
$ws.Range("B6").Value = 'PhidgetInterfaceKit 0/16/16'
$ws.Range("E6").Value = '16 Digital Inputs and 16 Digital Outputs'
$ws.Range("G6").Value = 'X'
$ws.Range("A7").Value = 1056
$ws.Range("B7").Value = 'PhidgetSpatial 3/3/3'
$ws.Range("E7").Value = 'Compass 3-Axis, Gyroscope 3-Axis, Accelerometer 3-Axis 5G'
$ws.Range("G7").Value = 'X'
$ws.Range("A8").Value = 1061
$ws.Range("B8").Value = 'PhidgetAdvancedServo 8-Motor'
$ws.Range("E8").Value = 'Contol the position, velocity, and acceleration of up to 8 RC servo motors'
$ws.Range("G8").Value = 'X'
$ws.Range("A9").Value = 1062
$ws.Range("B9").Value = 'PhidgetStepper Unipolar 4-Motor'
$ws.Range("E9").Value = 'Contol the position, velocity, and acceleration of up to 4 unipolar stepper motors'
$ws.Range("G9").Value = 'X'
$ws.Range("A10").Value = 1063
$ws.Range("B10").Value = 'PhidgetStepper Bipolar 1-Motor'
$ws.Range("E10").Value = 'Contol the position, velocity, and acceleration of 1 bipolar stepper motor'
$ws.Range("G10").Value = 'X'
$ws.Range("A11").Value = 1064
$ws.Range("B11").Value = 'PhidgetMotorControl HC'
$ws.Range("E11").Value = 'Control the velocity and acceleration of up to two high-current DC motors'
$ws.Range("G11").Value = 'X'
$ws.Range("A12").Value = 1070
$ws.Range("B12").Value = 'PhidgetSBC'
$ws.Range("C12").Value = 'psbc1, 192.168.150.1'
$ws.Range("E12").Value = 'Fully functional Single Board Computer with integrated 8/8/8 Interface Kit'
$ws.Range("G12").Value = 'X'
$ws.Range("A13").Value = 1072
$ws.Range("B13").Value = 'PhidgetSBC2'
$ws.Range("C13").Value = 'psbc21, 192.168.150.21'
$ws.Range("D13").Value = 48284
$ws.Range("E13").Value = 'Fully functional Single Board Computer with integrated 8/8/8 Interface Kit'
$ws.Range("G13").Value = 'X'
$ws.Range("A14").Value = 1072
$ws.Range("B14").Value = 'PhidgetSBC2'
$ws.Range("C14").Value = 'psbc22, 192.168.150.22'
$ws.Range("D14").Value = 48301
$ws.Range("E14").Value = 'Fully functional Single Board Computer with integrated 8/8/8 Interface Kit'
$ws.Range("G14").Value = 'X'
$ws.Range("A15").Value = 1072
$ws.Range("B15").Value = 'PhidgetSBC2'
$ws.Range("C15").Value = 'psbc23, 192.168.150.23'
$ws.Range("D15").Value = 251831
$ws.Range("E15").Value = 'Fully functional Single Board Computer with integrated 8/8/8 Interface Kit'
$ws.Range("G15").Value = 'X'
$ws.Range("A16").Value = 1101
$ws.Range("B16").Value = 'IR Distance Adapter'
$ws.Range("E16").Value = 'Bridge compatible Sharp IR sensors to PhidgetInterfaceKit'
$ws.Range("G16").Value = 'X'
$ws.Range("A17").Value = 1102
$ws.Range("B17").Value = 'IR Reflective Sensor 5mm'
$ws.Range("A18").Value = 1103
$ws.Range("B18").Value = 'IR Reflective Sensor 10cm'
$ws.Range("A19").Value = 1104
$ws.Range("B19").Value = 'Vibration Sensor'
$ws.Range("A20").Value = 1106
$ws.Range("B20").Value = 'Force Sensor'
$ws.Range("A21").Value = 1108
$ws.Range("B21").Value = 'Magnetic Sensor'
$ws.Range("A22").Value = 1109
$ws.Range("B22").Value = 'Rotation Sensor'
$ws.Range("A23").Value = 1110
$ws.Range("B23").Value = 'Touch Sensor'
$ws.Range("A24").Value = 1111
$ws.Range("B24").Value = 'Motion Sensor'
$ws.Range("E24").Value = 'Detects changes in infrared radiation when there is movement of an object different in temperature from the surroundings'
$ws.Range("G24").Value = 'X'
$ws.Range("A25").Value = 1112
$ws.Range("B25").Value = 'Slider Sensor'
$ws.Range("A26").Value = 1113
$ws.Range("B26").Value = 'Mini Joy Stick Sensor'
$ws.Range("A27").Value = 1124
$ws.Range("B27").Value = 'Precision Temperature Sensor'
$ws.Range("A28").Value = 1127
$ws.Range("B28").Value = 'Precision Light Sensor'
$ws.Range("A29").Value = 1128
$ws.Range("B29").Value = 'Sonar Sensor'
$ws.Range("E29").Value = 'Detects objects from 0 to 254 inches with 1 inch resolution'
$ws.Range("G29").Value = 'X'
$ws.Range("A30").Value = 3004
$ws.Range("B30").Value = 'Sensor Cables'
$ws.Range("A31").Value = 3051
$ws.Range("B31").Value = 'Dual Relay Board'

# ------------------------------------------------------------------
# 4. Formatting: wrap the Description column, center+wrap the
#    "Printed Manual" header only (data cells stay centered)
# ------------------------------------------------------------------
$ws.Range("G4").WrapText = $true
$ws.Range("E:E").WrapText = $true

# ------------------------------------------------------------------
# 5. Column widths / row heights
# ------------------------------------------------------------------
$ws.Range("A:A").ColumnWidth = 8.1796875
$ws.Range("B:B").ColumnWidth = 31.08984375
$ws.Range("C:C").ColumnWidth = 20.36328125
$ws.Range("E:E").ColumnWidth = 85.08984375
$ws.Range("G:G").ColumnWidth = 9.81640625
$ws.Range("H:H").ColumnWidth = 20.36328125

$ws.Rows.Item(4).RowHeight = 29
$ws.Rows.Item(24).RowHeight = 29

$ws.Cells.StandardHeight = 14.5

# ------------------------------------------------------------------
# 6. Rebuild Table1 over the new range so its column headers/ids are
#    sourced from the header row we just wrote
# ------------------------------------------------------------------
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A4:H31"), [System.Reflection.Missing]::Value, 1)
$newTbl.Name = "Table1"

# ------------------------------------------------------------------
# 7. Leave the selection where the author ended up
# ------------------------------------------------------------------
$ws.Range("I43").Select()
